$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1677
$ws1.Range("F5").Value = 756
$ws1.Range("F6").Value = 175

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1677
$ws4.Range("F6").Value = 756
$ws4.Range("F7").Value = 175
